$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wrapper_ready")
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
